$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates continue as serial numbers from existing data)
$newRows = @(
    @{ Row = 234; A = 44308; B = 1; C = 16; D = 162.2388967755019 },
    @{ Row = 235; A = 44309; B = 8; C = 22; D = 223.0784830663152 },
    @{ Row = 236; A = 44310; B = 2; C = 21; D = 212.9385520178463 },
    @{ Row = 237; A = 44311; B = 1; C = 19; D = 192.6586899209085 },
    @{ Row = 238; A = 44312; B = 0; C = 18; D = 182.5187588724397 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D

    # Copy the date-column formatting/style from the row above (column A),
    # matching the existing style used for all prior date cells.
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
